# Updates cryptocurrency Price (D) and Volume/1h (E) columns for rows 2-51
# Mirrors the source workbook's scheduled data refresh (GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'44.414.01"
$ws.Range("E2").Value = "  +0.45%  "
$ws.Range("D3").Value = "'2.224.29"
$ws.Range("E3").Value = "  -0.44%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'303.41"
$ws.Range("E5").Value = "  -0.68%  "
$ws.Range("D6").Value = "'90.39"
$ws.Range("E6").Value = "  -2.72%  "
$ws.Range("D7").Value = "'0.558"
$ws.Range("E7").Value = "  -2.00%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").Value = "'0.499"
$ws.Range("E9").Value = "  -3.66%  "
$ws.Range("D10").Value = "'33.88"
$ws.Range("E10").Value = "  +0.01%  "
$ws.Range("D11").Value = "'0.0783"
$ws.Range("E11").Value = "  -2.58%  "
$ws.Range("D12").Value = "'6.97"
$ws.Range("E12").Value = "  -1.51%  "
$ws.Range("E13").Value = "  -0.26%  "
$ws.Range("D14").Value = "'2.564.23"
$ws.Range("E14").Value = "  -0.41%  "
$ws.Range("D15").Value = "'2.267.19"
$ws.Range("E15").Value = "  -2.16%  "
$ws.Range("D16").Value = "'0.806"
$ws.Range("E16").Value = "  -1.77%  "
$ws.Range("D17").Value = "'13.21"
$ws.Range("E17").Value = "  -0.96%  "
$ws.Range("D18").Value = "'44.198.75"
$ws.Range("E18").Value = "  +0.65%  "
$ws.Range("D19").Value = "'0.0₃0911"
$ws.Range("E19").Value = "  -4.82%  "
$ws.Range("D20").Value = "'6.05"
$ws.Range("E20").Value = "  -4.02%  "
$ws.Range("D21").Value = "'11.40"
$ws.Range("E21").Value = "  -4.25%  "
$ws.Range("D22").Value = "'64.52"
$ws.Range("E22").Value = "  -1.35%  "
$ws.Range("D23").Value = "'234.11"
$ws.Range("E23").Value = "  -0.88%  "
$ws.Range("D24").Value = "'2.89"
$ws.Range("E24").Value = "  -6.13%  "
$ws.Range("E25").Value = "  -0.05%  "
$ws.Range("D26").Value = "'1.93"
$ws.Range("E26").Value = "  -2.84%  "
$ws.Range("D27").Value = "'2.27"
$ws.Range("E27").Value = "  +3.12%  "
$ws.Range("D28").Value = "'9.47"
$ws.Range("E28").Value = "  -3.19%  "
$ws.Range("D29").Value = "'36.40"
$ws.Range("E29").Value = "  -8.74%  "
$ws.Range("D30").Value = "'19.55"
$ws.Range("E30").Value = "  -1.88%  "
$ws.Range("D31").Value = "'5.65"
$ws.Range("E31").Value = "  -2.18%  "
$ws.Range("D32").Value = "'147.11"
$ws.Range("E32").Value = "  -2.80%  "
$ws.Range("D33").Value = "'2.62"
$ws.Range("E33").Value = "  +0.72%  "
$ws.Range("D34").Value = "'0.0758"
$ws.Range("E34").Value = "  -3.72%  "
$ws.Range("D35").Value = "'3.01"
$ws.Range("E35").Value = "  -1.00%  "
$ws.Range("E36").Value = "  -1.59%  "
$ws.Range("E37").Value = "  -2.64%  "
$ws.Range("D38").Value = "'1.79"
$ws.Range("E38").Value = "  +3.24%  "
$ws.Range("D39").Value = "'14.57"
$ws.Range("E39").Value = "  +2.70%  "
$ws.Range("D40").Value = "'3.24"
$ws.Range("E40").Value = "  -5.73%  "
$ws.Range("D41").Value = "'3.66"
$ws.Range("E41").Value = "  -2.28%  "
$ws.Range("D42").Value = "'0.0289"
$ws.Range("E42").Value = "  -2.02%  "
$ws.Range("E43").Value = "  -0.10%  "
$ws.Range("D44").Value = "'1.759.19"
$ws.Range("E44").Value = "  +3.30%  "
$ws.Range("E45").Value = "  +7.75%  "
$ws.Range("D46").Value = "'79.35"
$ws.Range("E46").Value = "  -2.75%  "
$ws.Range("D47").Value = "'0.182"
$ws.Range("E47").Value = "  -3.98%  "
$ws.Range("D48").Value = "'95.55"
$ws.Range("E48").Value = "  -3.11%  "
$ws.Range("D49").Value = "'4.73"
$ws.Range("E49").Value = "  -3.50%  "
$ws.Range("D50").Value = "'67.13"
$ws.Range("E50").Value = "  +1.60%  "
$ws.Range("D51").Value = "'52.51"
$ws.Range("E51").Value = "  -3.04%  "
